$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-99) holds a "Förändrad" (changed/modified) date that was
# bumped from 2023-10-05 (serial 45204) to 2023-10-08 (serial 45207) for
# every data row in the sheet.
for ($row = 2; $row -le 99; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
